$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Replace every single usage of old-style-12 cells (A6:A19, C14:C19 etc, F14:F19, G14:G18, H14, D19) with style 1 (donor B6)
$ws.Range("B6").Copy()
$ws.Range("A6:A19").PasteSpecial(-4122)
$ws.Range("C14:C19").PasteSpecial(-4122)
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("F14:F19").PasteSpecial(-4122)
$ws.Range("G14:G18").PasteSpecial(-4122)
$ws.Range("H14").PasteSpecial(-4122)
# Replace style 15 (H15:H18) with style 10 (donor H6)
$ws.Range("H6").Copy()
$ws.Range("H15:H18").PasteSpecial(-4122)
# Replace style 14 (A5) -> this is kept, no change needed since it just renumbers
